$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 currently holds the text "R40". It needs to become the text "1".
# A plain `$ws.Range("B11").Value = "1"` would be auto-coerced by Excel into
# the *number* 1 (since "1" parses as numeric), which would change the cell's
# type (t="s" -> no t attribute) and is not what the diff wants (it stays a
# shared-string cell, just pointing at a new shared string "1").
#
# Using a quote-prefixed literal ("'1") does force text, but it also flips
# the cell's style to a new "quote prefixed" style variant, which would
# change the cell's `s` attribute - again not what the diff shows (style
# stays s="23").
#
# So: compute "1" as text via a formula (TEXT() always yields a string),
# then paste only the *value* (not the format) over B11. That keeps B11's
# existing style/format untouched while making the stored cell a real text
# value.
$scratch = $ws.Range("Z100")
$scratch.Formula = '=TEXT(1,"0")'
$scratch.Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$scratch.Clear() | Out-Null
